$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (I1, J1) - copy formatting (style) from H1 which already has the header style
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows (no special style, like H2:H6)
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
